$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (spreadsheet row 3 since data starts at row 0)
$ws.Range("D2").Value = 1875587
$ws.Range("E2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiI0ODI5OGVhMC0yNDBhLTExZWUtOWMwNC1iMzcyMDk2MTViOGIiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiZGY2Y2UxNTEtYjQzMS00NjhhLWI1MTMtZjllNWZkOTdlMjMxIiwiaWF0IjoxNzAzODc2OTk5fQ.ga2AH4PHvUGlncTx7wiXm_DLWvUG98OVQR7wHOQMfjw"

# Add new row 3
$ws.Range("A3").Value = "RichDogeyBoy"
$ws.Range("B3").Value = "testact1112"
$ws.Range("C3").Value = 1500
$ws.Range("D3").Value = 5900
$ws.Range("E3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VySWQiOiIwOTQ4ZGE1MC04N2Q0LTExZWUtYjBjMi02MzM4M2I3OTUzNjAiLCJhcHBJZCI6MjMyLCJ0b2tlbklkIjoiNDE0YTdlYTAtMzE1My00YzZlLTk0MmItMjQ4N2FhZjdjNDQ1IiwiaWF0IjoxNzAzOTU0OTY5fQ.gYBD1eGIaV5ipOJPJUaAkH715hhkxMVPDwwn8GNddrY"
$ws.Range("F3").Value = "lwayzeyoihwk"
